$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values (columns B, C, D) for rows 4-11
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = 39.409999999999997
$ws.Range("D4").Value = 9.8000000000000007

$ws.Range("B5").Value = 101
$ws.Range("D5").Value = 55.020792079207922

$ws.Range("B6").Value = 108
$ws.Range("D6").Value = 62.776481481481497

$ws.Range("B7").Value = 375
$ws.Range("D7").Value = 41.899120000000032

$ws.Range("B8").Value = 297
$ws.Range("D8").Value = 47.807306397306391

$ws.Range("B9").Value = 221
$ws.Range("D9").Value = 38.557828054298618

$ws.Range("B10").Value = 118
$ws.Range("D10").Value = 40.057457627118623

$ws.Range("B11").Value = 105
$ws.Range("D11").Value = 39.846666666666657

# Update sheet view: zoom and selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 210
$ws.Range("G17").Select()
